$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to text format so numeric-looking strings
# (e.g. "224.05", "0.290") are not auto-converted to numbers by Excel,
# matching the inline-string cells in the source workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '33.799.55'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '1.779.58'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '224.05'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '32.17'
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('D9').Value = '0.290'
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('D10').Value = '0.0682'
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').Value = '2.033.11'
$ws.Range('D13').Value = '11.22'
$ws.Range('E13').Value = '  +4.52%  '
$ws.Range('D14').Value = '1.766.65'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '33.782.26'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '0.611'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('E17').Value = '  -2.32%  '
$ws.Range('D18').Value = '66.64'
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('D19').Value = '238.66'
$ws.Range('E19').Value = '  -2.90%  '
$ws.Range('D20').Value = '0.0₃0774'
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D22').Value = '10.60'
$ws.Range('E22').Value = '  -1.55%  '
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').Value = '160.42'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').Value = '16.12'
$ws.Range('E26').Value = '  -2.02%  '
$ws.Range('D27').Value = '7.03'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').Value = '0.0512'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').Value = '3.60'
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').Value = '3.50'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = '1.80'
$ws.Range('E34').Value = '  -2.35%  '
$ws.Range('D35').Value = '1.384.83'
$ws.Range('E35').Value = '  -2.04%  '
$ws.Range('D36').Value = '0.648'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('E39').Value = '  +5.71%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').Value = '0.911'
$ws.Range('E41').Value = '  -3.37%  '
$ws.Range('D42').Value = '78.34'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').Value = '13.61'
$ws.Range('E43').Value = '  +14.68%  '
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('E47').Value = '  +11.11%  '
$ws.Range('D48').Value = '107.68'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  -1.79%  '
$ws.Range('D50').Value = '1.933.48'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('E51').Value = '  +0.07%  '

# Restore default style on the touched range (drop the temporary text
# number format) so cell styling matches the original workbook.
$ws.Range("D2:E51").Style = "Normal"
